# BulkSiteTemplate.xlsx edit
#
# The commit changes the PARENT_SITE_ID value used in the two sample data
# rows (row 2 and row 3, column A) of Sheet1 from the numeric value
# 1003001 to the text value "1441001" (stored as a shared string, i.e.
# Excel cell type "s", rather than a plain numeric cell).
#
# Because "1441001" looks like a number, a plain `.Value = "1441001"`
# assignment would make Excel re-interpret it as a number again. Forcing
# the cell's number format to Text ("@") first makes Excel keep the
# assigned value as literal text, matching the target t="s" shared-string
# cell produced by the original edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A2:A3")
$rng.NumberFormat = "@"
$rng.Value = "1441001"
